## download articles with pandoc title blocks
##
## The first two paragraphs of the document held an italic book title
## ("On Pilgrimage") followed by a manual line-break + "Foreword" rule,
## then a bold "By Dorothy Day" byline paragraph. Replace them with a
## pandoc-style title block: a Title-styled paragraph with the
## publisher/year ("New York: Catholic Worker Books, 1948", split into
## one run per token, as pandoc's docx writer emits) followed by a plain
## "% Dorothy Day" author line.

$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)

# Range spanning both paragraphs (through p2's paragraph mark) so a single
# InsertXML call replaces both paragraphs' contents in one shot.
$targetRange = $d.Range($p1.Range.Start, $p2.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Title"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">New</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">York</w:t></w:r>
            <w:r><w:t xml:space="preserve">:</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Catholic</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Worker</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Books</w:t></w:r>
            <w:r><w:t xml:space="preserve">,</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">1948</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t xml:space="preserve">% Dorothy Day</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$targetRange.InsertXML($xml)
